$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1054.7727
$ws.Range("J17").Value = 1101.8422
$ws.Range("L17").Value = 3305.5266
$ws.Range("N17").Value = -3641.5266

$ws.Range("H52").Value = 1866.6666
$ws.Range("I52").Value = 300
$ws.Range("J52").Value = 5000
$ws.Range("K52").Value = 900
$ws.Range("L52").Value = 15000
$ws.Range("M52").Value = -740
$ws.Range("N52").Value = -15320

$ws.Range("H58").Value = 282.14285
$ws.Range("I58").Value = 226.92308
$ws.Range("J58").Value = 1000
$ws.Range("K58").Value = 680.76924
$ws.Range("L58").Value = 3000
$ws.Range("M58").Value = -530.76924
$ws.Range("N58").Value = -3300

$ws.Range("H74").Value = 3292.276
$ws.Range("I74").Value = 2538.6
$ws.Range("J74").Value = 3449.2917
$ws.Range("K74").Value = 2538.6
$ws.Range("L74").Value = 3449.2917
$ws.Range("M74").Value = -1602.6
$ws.Range("N74").Value = -5321.2917

$ws.Range("H76").Value = 4626.316
$ws.Range("I76").Value = 3069.4443
$ws.Range("J76").Value = 6027.5
$ws.Range("K76").Value = 3069.4443
$ws.Range("L76").Value = 6027.5
$ws.Range("M76").Value = -2754.4443
$ws.Range("N76").Value = -6657.5

$ws.Range("H77").Value = 3292.276
$ws.Range("I77").Value = 2538.6
$ws.Range("J77").Value = 3449.2917
$ws.Range("K77").Value = 12693
$ws.Range("L77").Value = 17246.4585
$ws.Range("M77").Value = -8013
$ws.Range("N77").Value = -26606.4585

$ws.Range("H79").Value = 4626.316
$ws.Range("I79").Value = 3069.4443
$ws.Range("J79").Value = 6027.5
$ws.Range("K79").Value = 3069.4443
$ws.Range("L79").Value = 6027.5
$ws.Range("M79").Value = -1977.4443
$ws.Range("N79").Value = -8211.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 1950
$ws.Range("I63").Value = 1950
$ws.Range("K63").Value = 1950
$ws.Range("M63").Value = -1264

$ws.Range("H66").Value = 1950
$ws.Range("I66").Value = 1950
$ws.Range("K66").Value = 9750
$ws.Range("M66").Value = -6318

$ws.Range("H88").Value = 2873.3572
$ws.Range("I88").Value = 1817.8
$ws.Range("J88").Value = 3459.7778
$ws.Range("K88").Value = 1817.8
$ws.Range("L88").Value = 3459.7778
$ws.Range("M88").Value = -1411.8
$ws.Range("N88").Value = -4271.7778

$ws.Range("H91").Value = 2873.3572
$ws.Range("I91").Value = 1817.8
$ws.Range("J91").Value = 3459.7778
$ws.Range("K91").Value = 1817.8
$ws.Range("L91").Value = 3459.7778
$ws.Range("M91").Value = -413.8
$ws.Range("N91").Value = -6267.7778

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1806.9474
$ws.Range("I105").Value = 1624.8
$ws.Range("J105").Value = 2490
$ws.Range("K105").Value = 1624.8
$ws.Range("L105").Value = 2490
$ws.Range("M105").Value = 122.2
$ws.Range("N105").Value = -5984

$ws.Range("H134").Value = 1427.2632
$ws.Range("I134").Value = 939.8
$ws.Range("J134").Value = 3255.25
$ws.Range("K134").Value = 2819.4
$ws.Range("L134").Value = 9765.75
$ws.Range("M134").Value = -284.3999999999996
$ws.Range("N134").Value = -14835.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 5981.25
$ws.Range("I62").Value = 7930
$ws.Range("J62").Value = 2733.3333
$ws.Range("K62").Value = 7930
$ws.Range("L62").Value = 2733.3333
$ws.Range("M62").Value = -7306
$ws.Range("N62").Value = -3981.3333

$ws.Range("H65").Value = 5981.25
$ws.Range("I65").Value = 7930
$ws.Range("J65").Value = 2733.3333
$ws.Range("K65").Value = 39650
$ws.Range("L65").Value = 13666.6665
$ws.Range("M65").Value = -36530
$ws.Range("N65").Value = -19906.6665

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 4655.5
$ws.Range("I3").Value = 3005.7144
$ws.Range("J3").Value = 9274.9
$ws.Range("K3").Value = 9017.143199999999
$ws.Range("L3").Value = 27824.7
$ws.Range("M3").Value = -8905.143199999999
$ws.Range("N3").Value = -28048.7

$ws.Range("H113").Value = 22443.412
$ws.Range("I113").Value = 384
$ws.Range("J113").Value = 31634.834
$ws.Range("K113").Value = 1152
$ws.Range("L113").Value = 94904.50199999999
$ws.Range("M113").Value = 1018
$ws.Range("N113").Value = -99244.50199999999

$ws.Range("H133").Value = 349019.8
$ws.Range("I133").Value = 3088.4211
$ws.Range("J133").Value = 1006289.5
$ws.Range("K133").Value = 9265.263300000001
$ws.Range("L133").Value = 3018868.5
$ws.Range("M133").Value = -4205.263300000001
$ws.Range("N133").Value = -3028988.5

$ws.Range("H140").Value = 33519.973
$ws.Range("I140").Value = 37454.84
$ws.Range("K140").Value = 112364.52
$ws.Range("M140").Value = -107184.52

$ws.Range("H141").Value = 43413.293
$ws.Range("I141").Value = 43413.293
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 130239.879
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = -125059.879
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 3933.3333
$ws.Range("I70").Value = 3942.8572
$ws.Range("J70").Value = 3900
$ws.Range("K70").Value = 3942.8572
$ws.Range("L70").Value = 3900
$ws.Range("M70").Value = -3672.8572
$ws.Range("N70").Value = -4440

$ws.Range("H73").Value = 3933.3333
$ws.Range("I73").Value = 3942.8572
$ws.Range("J73").Value = 3900
$ws.Range("K73").Value = 3942.8572
$ws.Range("L73").Value = 3900
$ws.Range("M73").Value = -3006.8572
$ws.Range("N73").Value = -5772

$ws.Range("H80").Value = 2700.7334
$ws.Range("I80").Value = 3502.5
$ws.Range("J80").Value = 2577.3845
$ws.Range("K80").Value = 3502.5
$ws.Range("L80").Value = 2577.3845
$ws.Range("M80").Value = -2504.5
$ws.Range("N80").Value = -4573.3845

$ws.Range("H83").Value = 2700.7334
$ws.Range("I83").Value = 3502.5
$ws.Range("J83").Value = 2577.3845
$ws.Range("K83").Value = 17512.5
$ws.Range("L83").Value = 12886.9225
$ws.Range("M83").Value = -12520.5
$ws.Range("N83").Value = -22870.9225

$ws.Range("H109").Value = 33499.875
$ws.Range("J109").Value = 33499.875
$ws.Range("L109").Value = 33499.875
$ws.Range("N109").Value = -35579.875

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1008.2414
$ws.Range("I122").Value = 893
$ws.Range("J122").Value = 1450
$ws.Range("K122").Value = 2679
$ws.Range("L122").Value = 4350
$ws.Range("M122").Value = -229
$ws.Range("N122").Value = -9250
